$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply header style (border/bold/centered) to the new O1:R1 cells by copying
# the format from the existing N1 header cell before writing new values.
$ws.Range("N1").Copy()
$ws.Range("O1:R1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Build the full A1:R25 data block (matching the "after" state of the sheet).
$data = New-Object 'object[,]' 25,18
$data[0,1] = 0
$data[0,2] = 1
$data[0,3] = 2
$data[0,4] = 3
$data[0,5] = 4
$data[0,6] = 5
$data[0,7] = 6
$data[0,8] = 7
$data[0,9] = 8
$data[0,10] = 9
$data[0,11] = 10
$data[0,12] = 11
$data[0,13] = 12
$data[0,14] = 13
$data[0,15] = 14
$data[0,16] = 15
$data[0,17] = 16
$data[1,0] = 0
$data[1,1] = 1.05
$data[1,2] = 1.002845661109816
$data[1,3] = 1.021610587078818
$data[1,4] = 1.008965102664857
$data[1,6] = 1
$data[1,8] = 1.042643369515094
$data[1,9] = 1.024941590803464
$data[1,10] = 1.032776843176386
$data[1,11] = 1.020300936482983
$data[1,13] = 1.012340543898715
$data[1,16] = 1.02
$data[1,17] = 1.034246285804805
$data[2,0] = 1
$data[2,1] = 1.05
$data[2,2] = 1.005944067331699
$data[2,3] = 1.023572653562669
$data[2,4] = 1.011338323847849
$data[2,6] = 1
$data[2,8] = 1.043207719264188
$data[2,9] = 1.026273992307919
$data[2,10] = 1.033912546913233
$data[2,11] = 1.021828155460988
$data[2,13] = 1.012783070248174
$data[2,16] = 1.02
$data[2,17] = 1.035046653672345
$data[3,0] = 2
$data[3,1] = 1.05
$data[3,2] = 1.007918023684863
$data[3,3] = 1.02482670691627
$data[3,4] = 1.012855986964398
$data[3,6] = 1
$data[3,8] = 1.043560027053661
$data[3,9] = 1.027121331606549
$data[3,10] = 1.03463445196783
$data[3,11] = 1.022801688714883
$data[3,13] = 1.013064542377111
$data[3,16] = 1.02
$data[3,17] = 1.035558008274588
$data[4,0] = 3
$data[4,1] = 1.05
$data[4,2] = 1.008743001838728
$data[4,3] = 1.025353866247952
$data[4,4] = 1.013491884739442
$data[4,6] = 1
$data[4,8] = 1.043706821466955
$data[4,9] = 1.027476361620783
$data[4,10] = 1.034938449457805
$data[4,11] = 1.023209569758779
$data[4,13] = 1.013182643146952
$data[4,16] = 1.02
$data[4,17] = 1.03578017339157
$data[5,0] = 4
$data[5,1] = 1.05
$data[5,2] = 1.008883955945498
$data[5,3] = 1.025446548247093
$data[5,4] = 1.013600938515012
$data[5,6] = 1
$data[5,8] = 1.043733398785171
$data[5,9] = 1.027538556574269
$data[5,10] = 1.034993647757169
$data[5,11] = 1.023280353705011
$data[5,13] = 1.013203520697329
$data[5,16] = 1.02
$data[5,17] = 1.035827997874844
$data[6,0] = 5
$data[6,1] = 1.05
$data[6,2] = 1.007936730888053
$data[6,3] = 1.02484549336495
$data[6,4] = 1.012871279646988
$data[6,6] = 1
$data[6,8] = 1.043567645395365
$data[6,9] = 1.02713355315188
$data[6,10] = 1.034650121032409
$data[6,11] = 1.022813852430627
$data[6,13] = 1.013069116965934
$data[6,16] = 1.02
$data[6,17] = 1.035589313294302
$data[7,0] = 6
$data[7,1] = 1.05
$data[7,2] = 1.003908766503074
$data[7,3] = 1.022291356003986
$data[7,4] = 1.009779314426756
$data[7,6] = 1
$data[7,8] = 1.042843727321151
$data[7,9] = 1.025404188166944
$data[7,10] = 1.033177610755478
$data[7,11] = 1.020828428260292
$data[7,13] = 1.012494834248263
$data[7,16] = 1.02
$data[7,17] = 1.03455270138234
$data[8,0] = 7
$data[8,1] = 1.05
$data[8,2] = 0.9965534934278805
$data[8,3] = 1.017649778393307
$data[8,4] = 1.004177980911952
$data[8,6] = 1
$data[8,8] = 1.04145917628474
$data[8,9] = 1.022228546517716
$data[8,10] = 1.030462846481931
$data[8,11] = 1.017203635036173
$data[8,13] = 1.011439825134393
$data[8,16] = 1.02
$data[8,17] = 1.032629975502428
$data[9,0] = 8
$data[9,1] = 1.05
$data[9,2] = 0.9914750035781238
$data[9,3] = 1.014471678302285
$data[9,4] = 1.000343786087798
$data[9,6] = 1
$data[9,8] = 1.04046843654041
$data[9,9] = 1.020030986029122
$data[9,10] = 1.028584467354861
$data[9,11] = 1.014706836143699
$data[9,13] = 1.010710287275536
$data[9,16] = 1.02
$data[9,17] = 1.031318739343008
$data[10,0] = 9
$data[10,1] = 1.05
$data[10,2] = 0.9892403316519348
$data[10,3] = 1.013087531222318
$data[10,4] = 0.9986660845558027
$data[10,6] = 1
$data[10,8] = 1.040029169130926
$data[10,9] = 1.0190676785532
$data[10,10] = 1.027766785472203
$data[10,11] = 1.013613200487205
$data[10,13] = 1.010391325656256
$data[10,16] = 1.02
$data[10,17] = 1.030773624401439
$data[11,0] = 10
$data[11,1] = 1.05
$data[11,2] = 0.9883980836129264
$data[11,3] = 1.012562290322541
$data[11,4] = 0.9980343248059257
$data[11,6] = 1
$data[11,8] = 1.039859671854697
$data[11,9] = 1.018701705019628
$data[11,10] = 1.027452710242218
$data[11,11] = 1.013199281193833
$data[11,13] = 1.010269755768444
$data[11,16] = 1.02
$data[11,17] = 1.030551560657921
$data[12,0] = 11
$data[12,1] = 1.05
$data[12,2] = 0.9885780509375831
$data[12,3] = 1.012673595235656
$data[12,4] = 0.9981691192350737
$data[12,6] = 1
$data[12,8] = 1.039895413740981
$data[12,9] = 1.018779378517534
$data[12,10] = 1.027518718854391
$data[12,11] = 1.013287328165383
$data[12,13] = 1.010295476710246
$data[12,16] = 1.02
$data[12,17] = 1.030595723842901
$data[13,0] = 12
$data[13,1] = 1.05
$data[13,2] = 0.989170420365614
$data[13,3] = 1.013043536195088
$data[13,4] = 0.9986135595991068
$data[13,6] = 1
$data[13,8] = 1.040014897864928
$data[13,9] = 1.019037076003188
$data[13,10] = 1.027740244465329
$data[13,11] = 1.013578672859327
$data[13,13] = 1.010381125611413
$data[13,16] = 1.02
$data[13,17] = 1.030753784573658
$data[14,0] = 13
$data[14,1] = 1.05
$data[14,2] = 0.9895364849965115
$data[14,3] = 1.013274032430058
$data[14,4] = 0.9988886534482059
$data[14,6] = 1
$data[14,8] = 1.040089630129141
$data[14,9] = 1.019197363993517
$data[14,10] = 1.027879325792216
$data[14,11] = 1.013759514261931
$data[14,13] = 1.010434559193848
$data[14,16] = 1.02
$data[14,17] = 1.030858027559457
$data[15,0] = 14
$data[15,1] = 1.05
$data[15,2] = 0.9916398534140269
$data[15,3] = 1.014589223074229
$data[15,4] = 1.000470096062728
$data[15,6] = 1
$data[15,8] = 1.040509883145976
$data[15,9] = 1.020111274922294
$data[15,10] = 1.028664048553014
$data[15,11] = 1.014794215160116
$data[15,13] = 1.010738240875129
$data[15,16] = 1.02
$data[15,17] = 1.031416049459875
$data[16,0] = 15
$data[16,1] = 1.05
$data[16,2] = 0.9929469041303247
$data[16,3] = 1.015408127989825
$data[16,4] = 1.001455008141171
$data[16,6] = 1
$data[16,8] = 1.040768841853565
$data[16,9] = 1.02067884686415
$data[16,10] = 1.029151284863283
$data[16,11] = 1.015437573800142
$data[16,13] = 1.010926854667736
$data[16,16] = 1.02
$data[16,17] = 1.031763152575101
$data[17,0] = 16
$data[17,1] = 1.05
$data[17,2] = 0.9936995267496449
$data[17,3] = 1.015875537335693
$data[17,4] = 1.002022198271843
$data[17,6] = 1
$data[17,8] = 1.040914315276181
$data[17,9] = 1.021002689806316
$data[17,10] = 1.029425751726089
$data[17,11] = 1.015806116996748
$data[17,13] = 1.011034067894051
$data[17,16] = 1.02
$data[17,17] = 1.031945505219733
$data[18,0] = 17
$data[18,1] = 1.05
$data[18,2] = 0.9939595901570911
$data[18,3] = 1.016040798254864
$data[18,4] = 1.002218870763105
$data[18,6] = 1
$data[18,8] = 1.04096665940688
$data[18,9] = 1.021116780907637
$data[18,10] = 1.029525183881157
$data[18,11] = 1.015935074775888
$data[18,13] = 1.011072162827782
$data[18,16] = 1.02
$data[18,17] = 1.032022239089627
$data[19,0] = 18
$data[19,1] = 1.05
$data[19,2] = 0.9928067039311933
$data[19,3] = 1.015319863820228
$data[19,4] = 1.001349225463263
$data[19,6] = 1
$data[19,8] = 1.040740926912039
$data[19,9] = 1.020617759183784
$data[19,10] = 1.029098576677374
$data[19,11] = 1.015368391784481
$data[19,13] = 1.010906521487779
$data[19,16] = 1.02
$data[19,17] = 1.031724541916115
$data[20,0] = 19
$data[20,1] = 1.05
$data[20,2] = 0.989001917925609
$data[20,3] = 1.012943230695864
$data[20,4] = 0.9984878809009744
$data[20,6] = 1
$data[20,8] = 1.039983919369233
$data[20,9] = 1.018966771467488
$data[20,10] = 1.027683469748557
$data[20,11] = 1.013497957582574
$data[20,13] = 1.010358208054003
$data[20,16] = 1.02
$data[20,17] = 1.030727234711792
$data[21,0] = 20
$data[21,1] = 1.05
$data[21,2] = 0.9865613906268966
$data[21,3] = 1.011418091022678
$data[21,4] = 0.9966589122151729
$data[21,6] = 1
$data[21,8] = 1.039487475853106
$data[21,9] = 1.01790296325615
$data[21,10] = 1.02676669684197
$data[21,11] = 1.012296848546647
$data[21,13] = 1.010004393664113
$data[21,16] = 1.02
$data[21,17] = 1.030065446040619
$data[22,0] = 21
$data[22,1] = 1.05
$data[22,2] = 0.9878522864974081
$data[22,3] = 1.012218348164729
$data[22,4] = 0.9976246829110123
$data[22,6] = 1
$data[22,8] = 1.039747179023508
$data[22,9] = 1.018462170449726
$data[22,10] = 1.027244284003641
$data[22,11] = 1.012929436279999
$data[22,13] = 1.010189837853467
$data[22,16] = 1.02
$data[22,17] = 1.030393437003073
$data[23,0] = 22
$data[23,1] = 1.05
$data[23,2] = 0.9928587169656252
$data[23,3] = 1.015342764247938
$data[23,4] = 1.001386970894558
$data[23,6] = 1
$data[23,8] = 1.040745328043906
$data[23,9] = 1.020634463323221
$data[23,10] = 1.029105696525941
$data[23,11] = 1.015389773470927
$data[23,13] = 1.010911219987078
$data[23,16] = 1.02
$data[23,17] = 1.031702017621446
$data[24,0] = 23
$data[24,1] = 1.05
$data[24,2] = 0.9984972898695352
$data[24,3] = 1.018882797934781
$data[24,4] = 1.005654163359087
$data[24,6] = 1
$data[24,8] = 1.041837906336793
$data[24,9] = 1.023075332251415
$data[24,10] = 1.031194887747193
$data[24,11] = 1.018165335266326
$data[24,13] = 1.011721963382824
$data[24,16] = 1.02
$data[24,17] = 1.033176323182449

$ws.Range("A1:R25").Value2 = $data
